$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Construction Progresses on Six Mid-Rise Developments in Astoria, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2026/01/construction-progresses-on-six-mid-rise-developments-in-astoria-queens.html"
$ws.Range("C2").Value = 'YIMBY recently photographed the progress of six more mid-rise residential and commercial buildings under construction in <a href="https://newyorkyimby.com/neighborhoods/astoria">Astoria</a>, Queens. The projects range from five to 13 stories and are located around the main 31st Street corridor and its elevated subway stations serving the N and W trains.'
$ws.Range("D2").Value = "2026-01-25T12:30:05+00:00"
$ws.Range("E2").Value = "Sun, 25 Jan 2026 12:30:05 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Astoria"
$ws.Range("H2").Value = ""
